# This script applies a weekly data update to the "Pepino ensalada" sheet.
# Two new rows (a fresh week's "Primera" and "Segunda" quality price quotes)
# are inserted immediately above the existing data block starting at row 344,
# pushing all the subsequent rows down by two positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at the top of the data block (rows 344-345),
# shifting the previous rows 344..459 down to 346..461.
$ws.Rows("344:345").Insert()

# New row 344: "Primera" quality quote for the new week (2023-03-22)
$ws.Cells.Item(344, 1).Value  = 1
$ws.Cells.Item(344, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(344, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(344, 4).Value  = 45007
$ws.Cells.Item(344, 5).Value  = 15
$ws.Cells.Item(344, 6).Value  = 100112043
$ws.Cells.Item(344, 7).Value  = "Pepino ensalada"
$ws.Cells.Item(344, 8).Value  = "Sin especificar"
$ws.Cells.Item(344, 9).Value  = "Primera"
$ws.Cells.Item(344, 10).Value = 130
$ws.Cells.Item(344, 11).Value = 3500
$ws.Cells.Item(344, 12).Value = 4000
$ws.Cells.Item(344, 13).Value = 3750
$ws.Cells.Item(344, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(344, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(344, 16).Value = 54
$ws.Cells.Item(344, 17).Value = 70
$ws.Cells.Item(344, 18).Value = "Hortaliza"

# New row 345: "Segunda" quality quote for the same new week (2023-03-22)
$ws.Cells.Item(345, 1).Value  = 1
$ws.Cells.Item(345, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(345, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(345, 4).Value  = 45007
$ws.Cells.Item(345, 5).Value  = 15
$ws.Cells.Item(345, 6).Value  = 100112043
$ws.Cells.Item(345, 7).Value  = "Pepino ensalada"
$ws.Cells.Item(345, 8).Value  = "Sin especificar"
$ws.Cells.Item(345, 9).Value  = "Segunda"
$ws.Cells.Item(345, 10).Value = 150
$ws.Cells.Item(345, 11).Value = 3000
$ws.Cells.Item(345, 12).Value = 3500
$ws.Cells.Item(345, 13).Value = 3250
$ws.Cells.Item(345, 14).Value = "`$/caja 100 unidades"
$ws.Cells.Item(345, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(345, 16).Value = 32
$ws.Cells.Item(345, 17).Value = 100
$ws.Cells.Item(345, 18).Value = "Hortaliza"
